$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1039
$ws.Range("I6").Value = 797.6875
$ws.Range("K6").Value = 2393.0625
$ws.Range("M6").Value = -2281.0625
# Row 9
$ws.Range("H9").Value = 549.8333
$ws.Range("I9").Value = 599.75
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 599.75
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = -430.75
$ws.Range("N9").Value = -788
# Row 17
$ws.Range("H17").Value = 808.13635
$ws.Range("J17").Value = 799
$ws.Range("L17").Value = 2397
$ws.Range("N17").Value = -2733
# Row 40
$ws.Range("H40").Value = 3262.5
$ws.Range("I40").Value = 1900
$ws.Range("J40").Value = 4625
$ws.Range("K40").Value = 1900
$ws.Range("L40").Value = 4625
$ws.Range("M40").Value = -1725
$ws.Range("N40").Value = -4975
# Row 48
$ws.Range("H48").Value = 5571.4287
$ws.Range("I48").Value = 5500
$ws.Range("K48").Value = 16500
$ws.Range("M48").Value = -16208
# Row 56
$ws.Range("H56").Value = 5571.4287
$ws.Range("I56").Value = 5500
$ws.Range("K56").Value = 16500
$ws.Range("M56").Value = -15966
# Row 100
$ws.Range("H100").Value = 5926.7
$ws.Range("I100").Value = 6196.3887
$ws.Range("J100").Value = 3499.5
$ws.Range("K100").Value = 6196.3887
$ws.Range("L100").Value = 3499.5
$ws.Range("M100").Value = -5655.3887
$ws.Range("N100").Value = -4581.5
# Row 113
$ws.Range("H113").Value = 6628.7144
$ws.Range("I113").Value = 6802.0835
$ws.Range("K113").Value = 6802.0835
$ws.Range("M113").Value = -3548.0835
# Row 129
$ws.Range("H129").Value = 1355.5
$ws.Range("I129").Value = 861.125
$ws.Range("K129").Value = 2583.375
$ws.Range("M129").Value = 2416.625
# Row 137
$ws.Range("H137").Value = 38665
$ws.Range("I137").Value = 31123.9
$ws.Range("J137").Value = 51233.5
$ws.Range("K137").Value = 93371.70000000001
$ws.Range("L137").Value = 153700.5
$ws.Range("M137").Value = -90821.70000000001
$ws.Range("N137").Value = -158800.5
# Row 138
$ws.Range("H138").Value = 23922.348
$ws.Range("I138").Value = 1686.7941
$ws.Range("J138").Value = 86923.086
$ws.Range("K138").Value = 5060.3823
$ws.Range("L138").Value = 260769.258
$ws.Range("M138").Value = 79.61769999999979
$ws.Range("N138").Value = -271049.258

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 4177.4
$ws.Range("I28").Value = 4177.4
$ws.Range("K28").Value = 4177.4
$ws.Range("M28").Value = -3985.4
# Row 32
$ws.Range("H32").Value = 20792.773
$ws.Range("I32").Value = 23276.021
$ws.Range("K32").Value = 23276.021
$ws.Range("M32").Value = -22989.021
# Row 41
$ws.Range("H41").Value = 2166.6667
$ws.Range("I41").Value = 2166.6667
$ws.Range("K41").Value = 2166.6667
$ws.Range("M41").Value = -1752.6667
# Row 61
$ws.Range("H61").Value = 5287.154
$ws.Range("I61").Value = 1174.4117
$ws.Range("J61").Value = 13055.667
$ws.Range("K61").Value = 1174.4117
$ws.Range("L61").Value = 13055.667
$ws.Range("M61").Value = -962.4117000000001
$ws.Range("N61").Value = -13479.667
# Row 92
$ws.Range("H92").Value = 9999
$ws.Range("J92").Value = 9999
$ws.Range("L92").Value = 9999
$ws.Range("N92").Value = -14991
# Row 97
$ws.Range("H97").Value = 1443.7693
$ws.Range("I97").Value = 1117.5
$ws.Range("J97").Value = 2177.875
$ws.Range("K97").Value = 1117.5
$ws.Range("L97").Value = 2177.875
$ws.Range("M97").Value = -621.5
$ws.Range("N97").Value = -3169.875
# Row 99
$ws.Range("H99").Value = 4177.4
$ws.Range("I99").Value = 4177.4
$ws.Range("K99").Value = 4177.4
$ws.Range("M99").Value = -1182.4
# Row 122
$ws.Range("H122").Value = 1544.1177
$ws.Range("I122").Value = 1338.2903
$ws.Range("J122").Value = 3671
$ws.Range("K122").Value = 4014.8709
$ws.Range("L122").Value = 11013
$ws.Range("M122").Value = -1564.8709
$ws.Range("N122").Value = -15913
# Row 134
$ws.Range("H134").Value = 77331.664
$ws.Range("J134").Value = 77331.664
$ws.Range("L134").Value = 77331.664
$ws.Range("N134").Value = -87471.664
# Row 136
$ws.Range("H136").Value = 5287.154
$ws.Range("I136").Value = 1174.4117
$ws.Range("J136").Value = 13055.667
$ws.Range("K136").Value = 3523.2351
$ws.Range("L136").Value = 39167.001
$ws.Range("M136").Value = -973.2351000000003
$ws.Range("N136").Value = -44267.001

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 2125
$ws.Range("I5").Value = 1250
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 1250
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1137
$ws.Range("N5").Value = -3226
# Row 134
$ws.Range("H134").Value = 3140.1333
$ws.Range("I134").Value = 2675.6667
$ws.Range("K134").Value = 8027.000100000001
$ws.Range("M134").Value = -5492.000100000001

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 103
$ws.Range("H103").Value = 12583
$ws.Range("I103").Value = 12583
$ws.Range("K103").Value = 12583
$ws.Range("M103").Value = -11411

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 550.6667
$ws.Range("I2").Value = 550.6667
$ws.Range("K2").Value = 3304.0002
$ws.Range("M2").Value = -3191.0002
# Row 4
$ws.Range("H4").Value = 107745020
$ws.Range("I4").Value = 144079340
$ws.Range("K4").Value = 432238020
$ws.Range("M4").Value = -432237908
# Row 7
$ws.Range("H7").Value = 230196.67
$ws.Range("I7").Value = 230196.67
$ws.Range("K7").Value = 690590.01
$ws.Range("M7").Value = -690478.01
# Row 55
$ws.Range("H55").Value = 5000
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
# Row 88
$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
# Row 91
$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
# Row 131
$ws.Range("H131").Value = 105866.39
$ws.Range("I131").Value = 330288.53
$ws.Range("J131").Value = 1670.3928
$ws.Range("K131").Value = 990865.5900000001
$ws.Range("L131").Value = 5011.178400000001
$ws.Range("M131").Value = -985825.5900000001
$ws.Range("N131").Value = -15091.1784

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1296.1818
$ws.Range("I97").Value = 1115.7059
$ws.Range("K97").Value = 1115.7059
$ws.Range("M97").Value = -619.7058999999999
# Row 122
$ws.Range("H122").Value = 4680.25
$ws.Range("I122").Value = 4563
$ws.Range("K122").Value = 13689
$ws.Range("M122").Value = -11239
# Row 131
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1253.45
$ws.Range("I82").Value = 1337.2307
$ws.Range("J82").Value = 1097.8572
$ws.Range("K82").Value = 1337.2307
$ws.Range("L82").Value = 1097.8572
$ws.Range("M82").Value = -976.2307000000001
$ws.Range("N82").Value = -1819.8572
# Row 85
$ws.Range("H85").Value = 1253.45
$ws.Range("I85").Value = 1337.2307
$ws.Range("J85").Value = 1097.8572
$ws.Range("K85").Value = 1337.2307
$ws.Range("L85").Value = 1097.8572
$ws.Range("M85").Value = -89.23070000000007
$ws.Range("N85").Value = -3593.8572
# Row 122
$ws.Range("H122").Value = 3263.4285
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 2968.8
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 8906.400000000001
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -13806.4

Write-Host "All cell updates applied."
